$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''25.819.88'
$ws.Range("D2").Style = "Normal"
$ws.Range("D3").Value = '''1.640.17'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '''  +0.52%  '
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = '''1.002'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '''  +0.01%  '
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = '''215.78'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '''  -0.06%  '
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = '''0.5064'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '''  -0.79%  '
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = '''1.003'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '''  +0.08%  '
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = '''0.2584'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '''  +0.60%  '
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = '''0.06433'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '''  +1.54%  '
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = '''20.45'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '''  +5.32%  '
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = '''0.07797'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '''  +0.23%  '
$ws.Range("E11").Style = "Normal"
$ws.Range("E12").Value = '''  +0.06%  '
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = '''1.644.50'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '''  +0.60%  '
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = '''1.866.88'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '''  +0.50%  '
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = '''0.5624'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '''  +2.29%  '
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = '''0.0₅7665'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '''  +0.42%  '
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = '''63.42'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '''  -0.59%  '
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = '''25.834.77'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '''  -0.37%  '
$ws.Range("E18").Style = "Normal"
$ws.Range("E19").Value = '''  +0.14%  '
$ws.Range("E19").Style = "Normal"
$ws.Range("B20").Value = 'BitcoinCash'
$ws.Range("C20").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D20").Value = '''193.35'
$ws.Range("D20").Style = "Normal"
$ws.Range("B21").Value = 'Uniswap'
$ws.Range("C21").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D21").Value = '''4.387'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '''  -0.59%  '
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = '''9.940'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '''  +0.94%  '
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = '''6.145'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '''  +1.98%  '
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = '''1.002'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '''  +0.16%  '
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = '''1.801'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '''  -4.67%  '
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = '''139.41'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '''  -1.80%  '
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = '''0.1235'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '''  -1.20%  '
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").Value = '''6.835'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '''  +1.30%  '
$ws.Range("E28").Style = "Normal"
$ws.Range("D29").Value = '''15.58'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '''  +0.34%  '
$ws.Range("E29").Style = "Normal"
$ws.Range("D31").Value = '''0.04972'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '''  +1.90%  '
$ws.Range("E31").Style = "Normal"
$ws.Range("D32").Value = '''3.292'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '''  +1.86%  '
$ws.Range("E32").Style = "Normal"
$ws.Range("D33").Value = '''3.251'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '''  +2.21%  '
$ws.Range("E33").Style = "Normal"
$ws.Range("D34").Value = '''1.571'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '''  +2.07%  '
$ws.Range("E34").Style = "Normal"
$ws.Range("D35").Value = '''2.385'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '''  +0.49%  '
$ws.Range("E35").Style = "Normal"
$ws.Range("D36").Value = '''0.9040'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '''  +0.75%  '
$ws.Range("E36").Style = "Normal"
$ws.Range("D37").Value = '''2.571'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '''  +1.20%  '
$ws.Range("E37").Style = "Normal"
$ws.Range("D38").Value = '''0.5567'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '''  +1.29%  '
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = '''1.132.93'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '''  +1.51%  '
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = '''0.01572'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '''  +0.93%  '
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = '''0.9964'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '''  -0.40%  '
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = '''5.486'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '''  -1.64%  '
$ws.Range("E42").Style = "Normal"
$ws.Range("E43").Value = '''  +0.78%  '
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = '''99.00'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '''  +1.50%  '
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = '''1.776.75'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '''  +0.55%  '
$ws.Range("E45").Style = "Normal"
$ws.Range("E46").Value = '''  -3.64%  '
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = '''55.64'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '''  +1.79%  '
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = '''0.4264'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '''  -3.94%  '
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = '''7.800'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '''  +3.36%  '
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = '''0.05032'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '''  -1.90%  '
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = '''0.9961'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '''  -0.64%  '
$ws.Range("E51").Style = "Normal"
